$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "quotations": a new quotation (Dilmatec / SKB8DXCYGE) was
# created and now sits at row 3, pushing the previously-existing
# rows down by one.
# ------------------------------------------------------------------
$wsQ = $wb.Worksheets.Item("quotations")
$wsQ.Rows("3:3").Insert()

$wsQ.Cells.Item(3,1).Value  = "MmE3ZThjNjktMTIyMy00MWQzLTk2N2QtNDc0Nzg4ZDJmYjAxOjU3MDE2"
$wsQ.Cells.Item(3,2).Value  = "SKB8DXCYGE"
$wsQ.Cells.Item(3,3).Value  = "Dilmatec"
$wsQ.Cells.Item(3,4).Value  = "Diagnostico e resoluçao problema em 02 equipamentos que nao estavam atingindo temperatura"
$wsQ.Cells.Item(3,5).Value  = ""
$wsQ.Cells.Item(3,6).Value  = $false
$wsQ.Cells.Item(3,7).NumberFormat = "@"
$wsQ.Cells.Item(3,7).Value  = "1535"
$wsQ.Cells.Item(3,8).NumberFormat = "@"
$wsQ.Cells.Item(3,8).Value  = "1535"
$wsQ.Cells.Item(3,9).Value  = "Pendente"
$wsQ.Cells.Item(3,10).Value = "2025-09-04T15:24:13.519Z"
$wsQ.Cells.Item(3,11).Value = ""
$wsQ.Cells.Item(3,12).Value = ""
$wsQ.Cells.Item(3,13).Value = "Adriana Vieira Masini"
$wsQ.Cells.Item(3,14).Value = ""
$wsQ.Cells.Item(3,15).Value = "2025-08-28T15:35:41.479Z"
$wsQ.Cells.Item(3,16).Value = "NjU3ZmQ5M2UtYjlmZC00NTdmLTlmM2EtZGI0YzE2MmZhMGQ3OjU3MDE2"
$wsQ.Cells.Item(3,17).Value = "percentage"
$wsQ.Cells.Item(3,18).NumberFormat = "@"
$wsQ.Cells.Item(3,18).Value = "0"
$wsQ.Cells.Item(3,19).NumberFormat = "@"
$wsQ.Cells.Item(3,19).Value = "0"
$wsQ.Cells.Item(3,20).Value = "NDIzOTk2OTo1NzAxNg=="
$wsQ.Cells.Item(3,21).Value = "pending"

# ------------------------------------------------------------------
# Sheet "items": the new quotation above brought two new line items
# with it; they land at rows 3-4, pushing the previously-existing
# item rows down by two.
# ------------------------------------------------------------------
$wsI = $wb.Worksheets.Item("items")
$wsI.Rows("3:4").Insert()

$wsI.Cells.Item(3,1).Value  = "MTc3NWUyYTAtZjAxNy00NWQwLTg2ZTMtYWFiZjYzMzZhOWUzOjU3MDE2"
$wsI.Cells.Item(3,2).Value  = 1
$wsI.Cells.Item(3,3).Value  = 680
$wsI.Cells.Item(3,4).Value  = "WRVT.00021 REALIZADO SERVIÇO LIMPEZA E CARGA DE GAS  R$680,00"
$wsI.Cells.Item(3,5).Value  = 3
$wsI.Cells.Item(3,6).Value  = "MmE3ZThjNjktMTIyMy00MWQzLTk2N2QtNDc0Nzg4ZDJmYjAxOjU3MDE2"
$wsI.Cells.Item(3,7).Value  = "NWIwZWFlMmYtYjVkOC00NTU0LTkzZmYtZGM2ZGIwM2E1ZmEwOjU3MDE2"
$wsI.Cells.Item(3,8).Value  = 680
$wsI.Cells.Item(3,9).Value  = "product"
$wsI.Cells.Item(3,10).Value = "MmE3ZThjNjktMTIyMy00MWQzLTk2N2QtNDc0Nzg4ZDJmYjAxOjU3MDE2"

$wsI.Cells.Item(4,1).Value  = "NGMwNTQyYWMtNGNjZC00NjljLThlZWItMDYxMjg2NzkzMmJhOjU3MDE2"
$wsI.Cells.Item(4,2).Value  = 1
$wsI.Cells.Item(4,3).Value  = 855
$wsI.Cells.Item(4,4).Value  = "WRVT.00020  REALIZADO RECUPERAÇAO DA ESTAÇAO MICRO MOTOR  E REALIZADO LIMPEZA  NO SISTEMA E CARGA DE GAS"
$wsI.Cells.Item(4,5).Value  = 3
$wsI.Cells.Item(4,6).Value  = "MmE3ZThjNjktMTIyMy00MWQzLTk2N2QtNDc0Nzg4ZDJmYjAxOjU3MDE2"
$wsI.Cells.Item(4,7).Value  = "MWY3MGI1MWUtZWEwMC00YWEyLTgzZTItNDgwYzc2NzE1OTJkOjU3MDE2"
$wsI.Cells.Item(4,8).Value  = 855
$wsI.Cells.Item(4,9).Value  = "service"
$wsI.Cells.Item(4,10).Value = "MmE3ZThjNjktMTIyMy00MWQzLTk2N2QtNDc0Nzg4ZDJmYjAxOjU3MDE2"

Write-Output "edit applied"
